$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F3 picks up the "4:1" score, styled like the other score cells (e.g. F6)
$ws.Range("F6").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("F3").Value = "4:1"

# New legend row below the grid
$ws.Range("D15").Value = "targer"
$ws.Range("E15").Value = "defensor"

# Leave the selection where the author left it
$ws.Range("F2").Select()
